# Auto-generated Excel COM-interop script
#
# This workbook ("Gilgamesh_Profits") is a Final Fantasy XIV crafting-leve
# profit tracker: one table per crafting profession (ALC/ARM/BSM/CRP/CUL/
# GSM/LTW/WVR), each row a leve turn-in, columns H-N holding market-price
# snapshots (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
# None of these cells are formulas -- they are plain cached values refreshed
# by a scheduled market-data runner. This script reproduces that refresh:
# per-cell value overwrites, one cell value cleared (no HQ price available),
# and one new cell populated (HQ profit now computable).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (32 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 491.42105
$ws.Range("I4").Value = 411.05884
$ws.Range("K4").Value = 411.05884
$ws.Range("M4").Value = -297.05884
$ws.Range("H15").Value = 3423.319
$ws.Range("I15").Value = 3423.319
$ws.Range("K15").Value = 10269.957
$ws.Range("M15").Value = -10100.957
$ws.Range("H40").Value = 6004.1
$ws.Range("I40").Value = 4340.6665
$ws.Range("K40").Value = 4340.6665
$ws.Range("M40").Value = -4165.6665
$ws.Range("H96").Value = 639.9
$ws.Range("I96").Value = 340.125
$ws.Range("K96").Value = 1020.375
$ws.Range("M96").Value = 352.625
$ws.Range("H111").Value = 5000
$ws.Range("I111").Value = 5000
$ws.Range("K111").Value = 15000
$ws.Range("M111").Value = -11933
$ws.Range("H121").Value = 1804
$ws.Range("J121").Value = 1804
$ws.Range("L121").Value = 5412
$ws.Range("N121").Value = -8906
$ws.Range("H125").Value = 1432.6666
$ws.Range("I125").Value = 1399.5
$ws.Range("K125").Value = 12595.5
$ws.Range("M125").Value = -10135.5
$ws.Range("H135").Value = 2650.1
$ws.Range("I135").Value = 2633.4443
$ws.Range("K135").Value = 23700.9987
$ws.Range("M135").Value = -21165.9987

# ---- Sheet: ARM (34 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 33036.57
$ws.Range("I45").Value = 48328.445
$ws.Range("J45").Value = 5511.2
$ws.Range("K45").Value = 48328.445
$ws.Range("L45").Value = 5511.2
$ws.Range("M45").Value = -47951.445
$ws.Range("N45").Value = -6265.2
$ws.Range("H61").Value = 3908.6667
$ws.Range("I61").Value = 3890.4
$ws.Range("K61").Value = 3890.4
$ws.Range("M61").Value = -3678.4
$ws.Range("H88").Value = 6034
$ws.Range("I88").Value = 3092
$ws.Range("K88").Value = 3092
$ws.Range("M88").Value = -2686
$ws.Range("H91").Value = 6034
$ws.Range("I91").Value = 3092
$ws.Range("K91").Value = 3092
$ws.Range("M91").Value = -1688
$ws.Range("H97").Value = 733.3714
$ws.Range("I97").Value = 689.03845
$ws.Range("J97").Value = 861.44446
$ws.Range("K97").Value = 689.03845
$ws.Range("L97").Value = 861.44446
$ws.Range("M97").Value = -193.03845
$ws.Range("N97").Value = -1853.44446
$ws.Range("H122").Value = 4977.1333
$ws.Range("I122").Value = 1472.1666
$ws.Range("K122").Value = 4416.4998
$ws.Range("M122").Value = -1966.4998
$ws.Range("H136").Value = 3908.6667
$ws.Range("I136").Value = 3890.4
$ws.Range("K136").Value = 11671.2
$ws.Range("M136").Value = -9121.200000000001

# ---- Sheet: BSM (26 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 35726244
$ws.Range("I20").Value = 45468170
$ws.Range("K20").Value = 45468170
$ws.Range("M20").Value = -45467923
$ws.Range("H86").Value = 1830.2273
$ws.Range("I86").Value = 1560.375
$ws.Range("J86").Value = 2549.8333
$ws.Range("K86").Value = 1560.375
$ws.Range("L86").Value = 2549.8333
$ws.Range("M86").Value = -437.375
$ws.Range("N86").Value = -4795.8333
$ws.Range("H89").Value = 1830.2273
$ws.Range("I89").Value = 1560.375
$ws.Range("J89").Value = 2549.8333
$ws.Range("K89").Value = 7801.875
$ws.Range("L89").Value = 12749.1665
$ws.Range("M89").Value = -2185.875
$ws.Range("N89").Value = -23981.1665
$ws.Range("H94").Value = 86958050
$ws.Range("I94").Value = 153846690
$ws.Range("K94").Value = 153846690
$ws.Range("M94").Value = -153846239
$ws.Range("H105").Value = 16252316
$ws.Range("I105").Value = 835004.2
$ws.Range("K105").Value = 835004.2
$ws.Range("M105").Value = -833257.2

# ---- Sheet: CRP (40 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1158
$ws.Range("J22").Value = 450
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -1150
$ws.Range("H31").Value = 3602.8333
$ws.Range("I31").Value = 2740.0667
$ws.Range("J31").Value = 7916.6665
$ws.Range("K31").Value = 2740.0667
$ws.Range("L31").Value = 7916.6665
$ws.Range("M31").Value = -2445.0667
$ws.Range("N31").Value = -8506.666499999999
$ws.Range("H34").Value = 3602.8333
$ws.Range("I34").Value = 2740.0667
$ws.Range("J34").Value = 7916.6665
$ws.Range("K34").Value = 2740.0667
$ws.Range("L34").Value = 7916.6665
$ws.Range("M34").Value = -2538.0667
$ws.Range("N34").Value = -8320.666499999999
$ws.Range("H58").Value = 2039.4445
$ws.Range("I58").Value = 1142.5
$ws.Range("J58").Value = 3833.3333
$ws.Range("K58").Value = 1142.5
$ws.Range("L58").Value = 3833.3333
$ws.Range("M58").Value = -939.5
$ws.Range("N58").Value = -4239.3333
$ws.Range("H103").Value = 4903
$ws.Range("I103").Value = 4903
$ws.Range("K103").Value = 4903
$ws.Range("M103").Value = -3731
$ws.Range("H136").Value = 2039.4445
$ws.Range("I136").Value = 1142.5
$ws.Range("J136").Value = 3833.3333
$ws.Range("K136").Value = 3427.5
$ws.Range("L136").Value = 11499.9999
$ws.Range("M136").Value = -877.5
$ws.Range("N136").Value = -16599.9999
$ws.Range("H139").Value = 55499.5
$ws.Range("J139").Value = 55499.5
$ws.Range("L139").Value = 55499.5
$ws.Range("N139").Value = -65779.5

# ---- Sheet: CUL (19 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1557.5333
$ws.Range("I2").Value = 17.714285
$ws.Range("K2").Value = 106.28571
$ws.Range("M2").Value = 6.714290000000005
$ws.Range("H38").Value = 425.53845
$ws.Range("I38").Value = 90.28570999999999
$ws.Range("J38").Value = 816.6667
$ws.Range("K38").Value = 270.85713
$ws.Range("L38").Value = 2450.0001
$ws.Range("M38").Value = 76.14287000000002
$ws.Range("N38").Value = -3144.0001
$ws.Range("H56").Value = 5876
$ws.Range("I56").Value = 5876
$ws.Range("K56").Value = 5876
$ws.Range("M56").Value = -5346
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# ---- Sheet: GSM (8 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5525.8125
$ws.Range("I122").Value = 3710.7144
$ws.Range("K122").Value = 11132.1432
$ws.Range("M122").Value = -8682.143199999999
$ws.Range("H132").Value = 3446.389
$ws.Range("J132").Value = 7000
$ws.Range("L132").Value = 21000
$ws.Range("N132").Value = -26060

# ---- Sheet: LTW (54 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8078.8667
$ws.Range("I7").Value = 7015.4165
$ws.Range("K7").Value = 7015.4165
$ws.Range("M7").Value = -6903.4165
$ws.Range("H16").Value = 831.93335
$ws.Range("I16").Value = 855.9231
$ws.Range("J16").Value = 676
$ws.Range("K16").Value = 855.9231
$ws.Range("L16").Value = 676
$ws.Range("M16").Value = -685.9231
$ws.Range("N16").Value = -1016
$ws.Range("H22").Value = 12635.637
$ws.Range("I22").Value = 26910.2
$ws.Range("K22").Value = 26910.2
$ws.Range("M22").Value = -26615.2
$ws.Range("H27").Value = 12635.637
$ws.Range("I27").Value = 26910.2
$ws.Range("K27").Value = 26910.2
$ws.Range("M27").Value = -26803.2
$ws.Range("H46").Value = 3272.2222
$ws.Range("I46").Value = 2993.75
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 2993.75
$ws.Range("L46").Value = 5500
$ws.Range("M46").Value = -2805.75
$ws.Range("N46").Value = -5876
$ws.Range("H55").Value = 685.4167
$ws.Range("I55").Value = 527.5
$ws.Range("K55").Value = 527.5
$ws.Range("M55").Value = -354.5
$ws.Range("H93").Value = 542.9167
$ws.Range("I93").Value = 542.9167
$ws.Range("K93").Value = 542.9167
$ws.Range("M93").Value = 705.0833
$ws.Range("H97").Value = 29344
$ws.Range("J97").Value = 29344
$ws.Range("L97").Value = 29344
$ws.Range("N97").Value = -31326
$ws.Range("H122").Value = 5467.44
$ws.Range("I122").Value = 4849.409
$ws.Range("K122").Value = 14548.227
$ws.Range("M122").Value = -12098.227
$ws.Range("H126").Value = 8078.8667
$ws.Range("I126").Value = 7015.4165
$ws.Range("K126").Value = 21046.2495
$ws.Range("M126").Value = -18576.2495
$ws.Range("H132").Value = 3753.6
$ws.Range("I132").Value = 2504.3333
$ws.Range("K132").Value = 7512.999899999999
$ws.Range("M132").Value = -4982.999899999999
$ws.Range("H136").Value = 4899.6665
$ws.Range("I136").Value = 5713.857
$ws.Range("K136").Value = 17141.571
$ws.Range("M136").Value = -14591.571

# ---- Sheet: WVR (23 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 50342
$ws.Range("I109").Value = 50342
$ws.Range("K109").Value = 50342
$ws.Range("M109").Value = -48955
$ws.Range("H122").Value = 62501884
$ws.Range("I122").Value = 2771
$ws.Range("J122").Value = 125001000
$ws.Range("K122").Value = 8313
$ws.Range("L122").Value = 375003000
$ws.Range("M122").Value = -5863
$ws.Range("N122").Value = -375007900
$ws.Range("H132").Value = 3957.5
$ws.Range("I132").Value = 3561.4375
$ws.Range("K132").Value = 10684.3125
$ws.Range("M132").Value = -8154.3125
$ws.Range("H136").Value = 2400
$ws.Range("I136").Value = 1925
$ws.Range("K136").Value = 5775
$ws.Range("M136").Value = -3225
$ws.Range("H140").Value = 94624.664
$ws.Range("J140").Value = 94624.664
$ws.Range("L140").Value = 94624.664
$ws.Range("N140").Value = -104984.664

